$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.603.46"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "2.234.73"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "269.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +14.09%  "
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.626"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0925"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +15.86%  "
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").Value = "2.572.45"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.14%  "
$ws.Range("D16").Value = "2.243.38"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("E17").Value = "  +2.48%  "
$ws.Range("D18").Value = "43.566.09"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.23%  "
$ws.Range("E28").Value = "  +5.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("E30").Value = "  +2.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0932"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.124"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("E36").Value = "  -4.91%  "
$ws.Range("E37").Value = "  -4.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +19.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.219"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("E48").Value = "  +3.10%  "
$ws.Range("E49").Value = "  +2.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.439"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("D51").Value = "2.458.06"
$ws.Range("E51").Value = "  +0.31%  "
